$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Myr Superion', ['{2}', 'Artifact Creature — Myr', 'Spend only mana produced by creatures to cast this spell.', '5/6'])"
$ws.Range("A3").Value = "('Phyrexian Metamorph', ['{3}{U/P}', 'Artifact Creature — Shapeshifter', '({U/P} can be paid with either {U} or 2 life.)', 'You may have Phyrexian Metamorph enter the battlefield as a copy of any artifact or creature on the battlefield, except it’s an artifact in addition to its other types.', '0/0'])"
$ws.Range("A4").Value = "('Priest of Urabrask', ['{2}{R}', 'Creature — Human Cleric', 'When Priest of Urabrask enters the battlefield, add {R}{R}{R}.', '2/1'])"
$ws.Range("A5").Value = "('Pristine Talisman', ['{3}', 'Artifact', '{T}: Add {C}. You gain 1 life.'])"
$ws.Range("A6").Value = "('Sheoldred, Whispering One', ['{5}{B}{B}', 'Legendary Creature — Praetor', 'Swampwalk (This creature can’t be blocked as long as defending player controls a Swamp.)', 'At the beginning of your upkeep, return target creature card from your graveyard to the battlefield.', 'At the beginning of each opponent’s upkeep, that player sacrifices a creature.', '6/6'])"
$ws.Range("A7").Value = "('Surgical Extraction', ['{B/P}', 'Instant', '({B/P} can be paid with either {B} or 2 life.)', 'Choose target card in a graveyard other than a basic land card. Search its owner’s graveyard, hand, and library for any number of cards with the same name as that card and exile them. Then that player shuffles their library.'])"

$ws.Range("A8:A33").EntireRow.Delete()
